# Added drawer holes for screw and magnets.
#
# Semantic changes applied to Sheet1 (the "variables" table, columns
# A=name, B=value/formula, D=unit/status label):
#   1. Row 52 "end_stop_spacing" is renamed to "switch_spacing" and gets a
#      unit label ("y") in column D.
#   2. A new row "frame_face_y" is inserted right before "frame_border"
#      (old row 57), computed as CEILING(frame_y/2 + wood_thickness*2, 1).
#   3. Two new rows "magnet_screw_spacing" and "magnet_screw_from_edge"
#      are inserted right after "magnet_z" (old row 59) and before
#      "frame_foot_z", with values 30 and 10 (unit label "new").
#   4. "magnet_z"'s unit label changes from "y" to "new".
#   5. "laser_y_offset" gets a unit label ("y") that it didn't have before.
#
# Excel auto-shifts every row-relative formula (B57->B58, B74->B77,
# B77->B80, ...) as a side effect of the real row inserts below, so we
# don't need to touch any of the untouched formulas by hand.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: new shared strings are appended to xl/sharedStrings.xml in the
# order they're first written into a cell, and the target workbook's
# table ends with ..., switch_spacing's old slot removed, then
# frame_face_y, magnet_screw_spacing, magnet_screw_from_edge,
# switch_spacing (in that order) appended at the end. So we add the
# brand-new labels (steps 1-2 below) before renaming
# end_stop_spacing -> switch_spacing (step 3), to reproduce that exact
# append order.

# ------------------------------------------------------------------
# 1. Insert "frame_face_y" before "frame_border" (currently row 57).
# ------------------------------------------------------------------
$ws.Rows("57:57").Insert()
$ws.Range("A57").Value = "frame_face_y"
$ws.Range("B57").Formula = "=CEILING(B56/2+B1*2,1)"
$ws.Range("D57").Value = "y"

# ------------------------------------------------------------------
# 2. Insert "magnet_screw_spacing" and "magnet_screw_from_edge" after
#    "magnet_z" (now row 60, since step 1 shifted everything by one).
# ------------------------------------------------------------------
$ws.Rows("61:62").Insert()

$ws.Range("A61").Value = "magnet_screw_spacing"
$ws.Range("B61").Value = 30
$ws.Range("D61").Value = "new"

$ws.Range("A62").Value = "magnet_screw_from_edge"
$ws.Range("B62").Value = 10
$ws.Range("D62").Value = "new"

# magnet_z (row 60) switches its unit label from "y" to "new".
$ws.Range("D60").Value = "new"

# ------------------------------------------------------------------
# 3. Rename end_stop_spacing -> switch_spacing (row 52) and label it.
# ------------------------------------------------------------------
$ws.Range("A52").Value = "switch_spacing"
$ws.Range("D52").Value = "y"

# ------------------------------------------------------------------
# 4. laser_y_offset (now row 65) gains a unit label it didn't have.
# ------------------------------------------------------------------
$ws.Range("D65").Value = "y"

# ------------------------------------------------------------------
# 5. Restore the view: selection on D53 (matches the committed state).
# ------------------------------------------------------------------
$ws.Range("D53").Select()
